$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Diana%Aguiar de Sousa%NULL%0, H Bart%van der Worp%NULL%2, H Bart%van der Worp%NULL%0, Valeria%Caso%NULL%1, Charlotte%Cordonnier%NULL%1, Daniel%Strbian%NULL%2, Daniel%Strbian%NULL%0, George%Ntaios%NULL%1, Peter D%Schellinger%NULL%1, Else Charlotte%Sandset%NULL%1, NULL%NULL%NULL%0, NULL%NULL%NULL%0]'
$ws.Range("I2").Value = ''
$ws.Range("J2").Value = 'SAGE Publications'

$ws.Range("D3").Value = 'Routine care for chronic disease is an ongoing major challenge.
 We aimed to evaluate the global impact of COVID-19 on routine care for chronic diseases.
 An online survey was posted 31 March to 23 April 2020 targeted at healthcare professionals.
 202 from 47 countries responded.
 Most reported change in routine care to virtual communication.
 Diabetes, chronic obstructive pulmonary disease, and hypertension were the most impacted conditions due to reduction in access to care.
 80% reported the mental health of their patients worsened during COVID-19. It is important routine care continues in spite of the pandemic, to avoid a rise in non-COVID-19-related morbidity and mortality.
'
$ws.Range("E3").Value = '[Yogini V.%Chudasama%NULL%0, Clare L.%Gillies%NULL%1, Francesco%Zaccardi%NULL%0, Briana%Coles%NULL%1, Melanie J.%Davies%NULL%1, Samuel%Seidu%NULL%1, Kamlesh%Khunti%NULL%0]'
$ws.Range("I3").Value = ''
$ws.Range("J3").Value = 'Diabetes India. Published by Elsevier Ltd.'

$ws.Range("D4").Value = 'Background
id="Par1">Anxiety, depression and reduction of quality of life (QoL) are common in people with multiple sclerosis (pwMS).

 Fear of getting sick from COVID-19, government’s lockdown and the imposed social distancing might have had an impact on psychological distress and QoL.


Objectives
id="Par2">The aim of our study was to investigate anxiety, depression and QoL changes in pwMS during SARS-CoV-2 outbreak and lockdown in Italy.


Methods
id="Par3">67 pwMS with a previous (less than 6 months) neuropsychological evaluation before SARS-CoV-2 outbreak (T0) were re-evaluated at the time of the outbreak and lockdown in Italy (T1).

 They underwent a clinical and neurological evaluation and completed the State-Trait Anxiety Inventory (STAI-Y1), the Beck Depression Inventory second edition (BDI-II), and Multiple Sclerosis Quality of Life-54 (MsQoL-54) at T0 and T1. Benjamini–Hochberg procedure was applied to control the false discovery rate.


Results
id="Par4">BDI-II and STAI-Y1 scores did not change between T0 and T1. At T1, MsQoL-54 scores were higher on the satisfaction with sexual life and the social function subscales, and lower on the limitation due to emotional problems subscale.


Conclusions
id="Par5">This is the first study that evaluated mood and QoL levels before and during the lockdown due to COVID-19 pandemic in pwMS.

 No worsening of anxiety and depression levels was found.

 Contrariwise some improvements were noted on QoL, the most reliable regarding the sexual satisfaction and the social function.


'
$ws.Range("E4").Value = '[Rocco%Capuano%NULL%0, Manuela%Altieri%NULL%1, Alvino%Bisecco%NULL%1, Alessandro%d’Ambrosio%NULL%1, Renato%Docimo%NULL%1, Daniela%Buonanno%NULL%1, Federica%Matrone%NULL%1, Federica%Giuliano%NULL%1, Gioacchino%Tedeschi%NULL%1, Gabriella%Santangelo%NULL%1, Antonio%Gallo%antonio.gallo@unicampania.it%1]'
$ws.Range("I4").Value = ''
$ws.Range("J4").Value = 'Springer Berlin Heidelberg'

$ws.Range("E5").Value = '[Jing%Zhao%NULL%0, Hang%Li%NULL%1, David%Kung%NULL%1, Marc%Fisher%NULL%1, Ying%Shen%NULL%1, Renyu%Liu%NULL%1]'
$ws.Range("I5").Value = ''
$ws.Range("J5").Value = 'Lippincott Williams &amp; Wilkins'

$ws.Range("D6").Value = 'Objective
id="Par1">To report the understanding and decision-making of neuroimmunologists and their treatment of patients with multiple sclerosis (MS) during the early stages of the SARS-CoV-2 (COVID-19) outbreak.


Methods
id="Par2">A survey instrument was designed and distributed online to neurologists in April 2020.
Results
id="Par3">There were 250 respondents (response rate 21.8%).

 243 saw &gt;  = 10 MS patients in the prior 6 months (average 197 patients) and were analyzed further (92% USA, 8% Canada; average practice duration 16 years; 5% rural, 17% small city, 38% large city, 40% highly urbanized).

 Patient volume dropped an average of 79% (53–11 per month).

 23% were aware of patients self-discontinuing a DMT due to fear of COVID-19 with 43% estimated to be doing so against medical advice.

 65% of respondents reported deferring &gt;  = 1 doses of a DMT (49%), changing the dosing interval (34%), changing to home infusions (20%), switching a DMT (9%), and discontinuing DMTs altogether (8%) as a result of COVID-19. Changes in DMTs were most common with the high-efficacy therapies alemtuzumab, cladribine, ocrelizumab, rituximab, and natalizumab.

 35% made no changes to DMT prescribing.

 98% expressed worry about their patients contracting COVID-19 and 78% expressed the same degree of worry about themselves.

  &gt; 50% believed high-efficacy DMTs prolong viral shedding of SARS-CoV-2 and that B-cell therapies might prevent protective vaccine effects.

 Accelerated pace of telemedicine and practice model changes were identified as major shifts in practice.


Conclusions
id="Par4">Reported prescribing changes and practice disruptions due to COVID-19 may be temporary but could have a lasting influence on MS care.


Electronic supplementary material
The online version of this article (10.1007/s00415-020-10045-9) contains supplementary material, which is available to authorized users.


'
$ws.Range("E6").Value = '[Farrah J.%Mateen%fmateen@mgh.harvard.edu%0, Shawheen%Rezaei%NULL%2, Shawheen%Rezaei%NULL%0, Nicholas%Alakel%NULL%1, Brittany%Gazdag%NULL%1, Aditya Ravi%Kumar%NULL%1, Andre%Vogel%NULL%1]'
$ws.Range("I6").Value = ''
$ws.Range("J6").Value = 'Springer Berlin Heidelberg'

$ws.Range("E7").Value = '[Ljiljana%Radulovic%NULL%0, Jevto%Erakovic%NULL%1, Milovan%Roganovic%NULL%1]'
$ws.Range("I7").Value = ''
$ws.Range("J7").Value = 'Elsevier B.V.'

$ws.Range("D8").Value = '
              •
              RRMS patients had different degrees of fear of COVID-19 disease.
'
$ws.Range("E8").Value = '[Aleksandar%Stojanov%NULL%0, Marina%Malobabic%NULL%1, Vuk%Milosevic%NULL%1, Jelena%Stojanov%NULL%1, Slobodan%Vojinovic%NULL%1, Goran%Stanojevic%NULL%1, Milos%Stevic%NULL%1]'
$ws.Range("I8").Value = ''
$ws.Range("J8").Value = 'Elsevier B.V.'

$ws.Range("D9").Value = 'Objective: Neurological sequelae of SARS-CoV-2 infection have already been reported, but there is insufficient data about the impact of the pandemic on the management of the patients with chronic neurological diseases.
 We aim to analyze the effect of COVID-19 pandemic and social restriction rules on these fragile patients.
'
$ws.Range("E9").Value = '[Carla%Piano%NULL%0, Enrico%Di Stasio%NULL%1, Guido%Primiano%NULL%1, Delfina%Janiri%NULL%1, Marco%Luigetti%NULL%1, Giovanni%Frisullo%NULL%1, Catello%Vollono%NULL%1, Matteo%Lucchini%NULL%1, Valerio%Brunetti%NULL%1, Mauro%Monforte%NULL%1, Valeria%Guglielmi%NULL%1, Giacomo%Della Marca%NULL%1, Amelia%Evoli%NULL%1, Camillo%Marra%NULL%1, Massimiliano%Mirabella%NULL%1, Davide%Quaranta%NULL%1, Enzo%Ricci%NULL%1, Serenella%Servidei%NULL%1, Gabriella%Silvestri%NULL%1, Simone%Bellavia%NULL%1, Sara%Bortolani%NULL%1, Francesco%Bove%NULL%1, Riccardo%Di Iorio%NULL%1, Andrea%Di Paolantonio%NULL%1, Danilo%Genovese%NULL%1, Tamara%Ialongo%NULL%1, Maria Rita%Lo Monaco%NULL%1, Jessica%Marotta%NULL%1, Agata Katia%Patanella%NULL%1, Alessia%Perna%NULL%1, Martina%Petracca%NULL%1, Giorgia%Presicce%NULL%1, Vittorio%Riso%NULL%1, Eleonora%Rollo%NULL%1, Angela%Romano%NULL%1, Marina%Romozzi%NULL%1, Cristina%Sancricca%NULL%1, Irene%Scala%NULL%1, Gregorio%Spagni%NULL%1, Marcella%Solito%NULL%1, Luca%Tricoli%NULL%1, Paola%Zinzi%NULL%1, Paolo%Calabresi%NULL%1, Anna Rita%Bentivoglio%NULL%1]'
$ws.Range("I9").Value = ''
$ws.Range("J9").Value = 'Frontiers Media S.A.'

$ws.Range("C10").Value = 'Unknown Title'
$ws.Range("D10").Value = 'Unknown Abstract'
$ws.Range("E10").Value = '[]'
$ws.Range("F10").Value = 'not found'
$ws.Range("G10").Value = 'N/A'
$ws.Range("I10").Value = ''

$ws.Range("D11").Value = 'Background
id="Par1">Quarantine was the measure taken by governments to control the rapid spread of COVID-19. This restriction resulted in a sudden change in people’s lifestyle, leading to an increase in sedentary behavior and a related decrease in the practice of physical activity (PA).

 However, in neuromuscular diseases patients need to perform regular PA to counteract the negative consequences of the disease.

 Hence, the aim of this study was to estimate the levels of PA, measured as energy expenditure (MET–minute/week), among patients with neuromuscular disease (NMD) before and during the last week of quarantine.


Methods
id="Par2">A total of 268 Italian subjects, living in Sicily, completed an adapted version of the IPAQ-SF.

 Participants comprised 149 NMD, enrolled at the Neuromuscular Clinic of Palermo and 119 healthy subjects (control group).

 The SF-12 questionnaire was also administered to NMD.

 The Mann–Whitney U and the Kruskal–Wallis rank-sum tests were used for statistical analyses.


Results
id="Par3">We observed a significant decrease of the total weekly PA level during COVID-19 quarantine in both patients and controls.

 Moreover, a significant difference in the total weekly PA level was found depending on the presence of neuromuscular disease, impaired walking, gender and BMI.

 Finally, we found a correlation between SF-12 scores and the entity of the reduction of PA level during quarantine, thus confirming a relevant association with the quality of life in NMD.


Conclusion
id="Par4">Our study confirmed that COVID-19 quarantine has affected the practice of PA among both NMD and healthy controls.


'
$ws.Range("E11").Value = '[Vincenzo%Di Stefano%vincenzo19689@gmail.com%0, Giuseppe%Battaglia%NULL%2, Giuseppe%Battaglia%NULL%0, Valerio%Giustino%NULL%1, Andrea%Gagliardo%NULL%1, Michele%D’Aleo%NULL%1, Ottavio%Giannini%NULL%1, Antonio%Palma%NULL%1, Filippo%Brighina%NULL%1]'
$ws.Range("I11").Value = ''
$ws.Range("J11").Value = 'Springer Berlin Heidelberg'

$ws.Range("E12").Value = '[Guillermo F%López-Sánchez%gfls@um.es%0, Rubén%López-Bueno%rlopezbu@unizar.es%1, Alejandro%Gil-Salmerón%alejandro.gil.salmeron@uv.es%1, Roksana%Zauder%NULL%1, Maria%Skalska%NULL%1, Joanna%Jastrzębska%NULL%1, Zbigniew%Jastrzębski%NULL%1, Felipe B%Schuch%NULL%1, Igor%Grabovac%NULL%0, Mark A%Tully%NULL%1, Lee%Smith%NULL%0]'
$ws.Range("I12").Value = ''
$ws.Range("J12").Value = 'Oxford University Press'
